# Generate Report for Handoff
# Regenerates handoff info: new source-file GUID, refreshed handoff timestamps,
# and clears the (not-yet-populated) "Latest Target File"/"Latest Handback File"/
# "Latest Handback DateTime" columns for both locales.

$wb = $excel.ActiveWorkbook

$oldGuid = "d7caa322-5ac1-430e-ba12-0ef535f71a30"
$newGuid = "87284e10-ec28-4220-a9c3-71d732d68a1c"
$oldHash = "5a7239e98103a6ce42c8d111a00091670c71668c"
$newHash = "3f917de7507826885a17b7286a4ca115f903f025"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-18 16:58:34"

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Address() -eq "`$B`$2") {
        $h.TextToDisplay = "e2e\$newGuid.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-18 16:58:29"
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Address() -eq "`$A`$2") {
        $h.TextToDisplay = "$newGuid.md"
    }
}

# Remove the "Latest Target File" hyperlink/value (I2) - not yet handed back
foreach ($h in @($wsZhCn.Hyperlinks)) {
    if ($h.Range.Address() -eq "`$I`$2") {
        $h.Delete()
    }
}
$wsZhCn.Range("I2").ClearContents()
$wsZhCn.Range("I2").Style = "Normal"

# Clear "Latest Handback File" (J2) - not yet handed back
$wsZhCn.Range("J2").ClearContents()
$wsZhCn.Range("J2").Style = "Normal"

$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-18 16:58:34"
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Address() -eq "`$A`$2") {
        $h.TextToDisplay = "$newGuid.md"
    }
}

# Remove the "Latest Target File" hyperlink/value (I2) - not yet handed back
foreach ($h in @($wsDeDe.Hyperlinks)) {
    if ($h.Range.Address() -eq "`$I`$2") {
        $h.Delete()
    }
}
$wsDeDe.Range("I2").ClearContents()
$wsDeDe.Range("I2").Style = "Normal"

# Clear "Latest Handback File" (J2) - not yet handed back
$wsDeDe.Range("J2").ClearContents()
$wsDeDe.Range("J2").Style = "Normal"

$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
